$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.510.39'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '1.691.06'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.52'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3894'
$ws.Range('E7').Value = '  -1.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4035'
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.003'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.58'
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08761'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.20'
$ws.Range('E13').Value = '  +6.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.530'
$ws.Range('E14').Value = '  +3.56%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.986'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001352'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = '1.685.01'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '98.47'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07107'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.05'
$ws.Range('E20').Value = '  +2.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.281'
$ws.Range('E21').Value = '  +4.03%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.29'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').Value = '24.497.69'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.976'
$ws.Range('E25').Value = '  -9.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.350'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.77'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.54'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.795'
$ws.Range('E29').Value = '  +17.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '137.21'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.227'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').Value = '1.871.38'
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.08847'
$ws.Range('E33').Value = '  +2.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.454'
$ws.Range('E34').Value = '  +5.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.039'
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2815'
$ws.Range('E36').Value = '  +2.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.960'
$ws.Range('E37').Value = '  +4.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02921'
$ws.Range('E38').Value = '  +7.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.81'
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09130'
$ws.Range('E41').Value = '  -1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7937'
$ws.Range('E42').Value = '  +3.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.456'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.64'
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.621'
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7255'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.205'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'Flow'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.354'
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.56'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '91.06'
$ws.Range('E51').Value = '  +0.08%  '
